$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "2003年" row (row 2) and "2008年" row (row 3)
# so remaining rows (2013年, 2018年) shift up to rows 2 and 3.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
